# daily auto push: 2026-01-19 18:48 UTC
#
# Two new 15-minute-interval log entries were appended for the current day
# ("2026/01/19" 22:00 and "2026/01/20" 02:00), pushing the whole tail of the
# rolling log (formerly rows 672-713, the "2026/12/29" .. "2027/01/05" block)
# down by two rows (now rows 674-715). The sheet's used range grows from
# A1:D713 to A1:D715.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the existing tail of the table down by inserting two blank rows
# right before the old row 672 (Excel renumbers/moves everything below
# automatically, and the sheet dimension is recalculated).
$ws.Rows("672:673").Insert()

# --- New row 672: 2026/01/19, 月, 22:00, rank 201 -------------------------
# Column A holds the date as literal text (matching the rest of the sheet,
# which stores dates as plain strings rather than real date serials), so we
# briefly force Text format before writing it and then drop the format again
# so the cell is left with the sheet's normal (default) style.
$ws.Cells.Item(672, 1).NumberFormat = "@"
$ws.Cells.Item(672, 1).Value = "2026/01/19"
$ws.Cells.Item(672, 1).ClearFormats()
$ws.Cells.Item(672, 2).Value = "月"
$ws.Cells.Item(672, 3).Value = 22
$ws.Cells.Item(672, 4).Value = 201

# --- New row 673: 2026/01/20, 火, 02:00, rank 201 -------------------------
$ws.Cells.Item(673, 1).NumberFormat = "@"
$ws.Cells.Item(673, 1).Value = "2026/01/20"
$ws.Cells.Item(673, 1).ClearFormats()
$ws.Cells.Item(673, 2).Value = "火"
$ws.Cells.Item(673, 3).Value = 2
$ws.Cells.Item(673, 4).Value = 201
